$d = $word.ActiveDocument

# Replace lowercase "historical" -> "recent"
$d.Content.Find.Execute("historical", $true, $true, $false, $false, $false,
                         $true, 1, $false, "recent", 2)

# Replace capitalized "Historical" -> "Recent"
$d.Content.Find.Execute("Historical", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Recent", 2)
